$wb = $excel.ActiveWorkbook

# ALC row 8
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 669.8182
$ws.Range("I8").Value = 40.88889
$ws.Range("J8").Value = 3500
$ws.Range("K8").Value = 122.66667
$ws.Range("L8").Value = 10500
$ws.Range("M8").Value = 16.33332999999999
$ws.Range("N8").Value = -10778

# ALC row 10
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").Value = ""

# ALC row 70
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 25005250
$ws.Range("I70").Value = 9000
$ws.Range("J70").Value = 33337334
$ws.Range("K70").Value = 27000
$ws.Range("L70").Value = 100012002
$ws.Range("M70").Value = -26730
$ws.Range("N70").Value = -100012542

# ALC row 73
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 25005250
$ws.Range("I73").Value = 9000
$ws.Range("J73").Value = 33337334
$ws.Range("K73").Value = 27000
$ws.Range("L73").Value = 100012002
$ws.Range("M73").Value = -26064
$ws.Range("N73").Value = -100013874

# ALC row 75
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 17657
$ws.Range("J75").Value = 17657
$ws.Range("L75").Value = 17657
$ws.Range("N75").Value = -19529

# ALC row 78
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 17657
$ws.Range("J78").Value = 17657
$ws.Range("L78").Value = 52971
$ws.Range("N78").Value = -62331

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2191.7021
$ws.Range("I137").Value = 1519.7297
$ws.Range("J137").Value = 4678
$ws.Range("K137").Value = 4559.189100000001
$ws.Range("L137").Value = 14034
$ws.Range("M137").Value = -2009.189100000001
$ws.Range("N137").Value = -19134

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 750.5
$ws.Range("I2").Value = 663
$ws.Range("J2").Value = 1013
$ws.Range("K2").Value = 663
$ws.Range("L2").Value = 1013
$ws.Range("M2").Value = -550
$ws.Range("N2").Value = -1239

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13724.018
$ws.Range("I32").Value = 10218.725
$ws.Range("J32").Value = 17229.31
$ws.Range("K32").Value = 10218.725
$ws.Range("L32").Value = 17229.31
$ws.Range("M32").Value = -9931.725
$ws.Range("N32").Value = -17803.31

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 750.5
$ws.Range("I116").Value = 663
$ws.Range("J116").Value = 1013
$ws.Range("K116").Value = 663
$ws.Range("L116").Value = 1013
$ws.Range("M116").Value = 1631
$ws.Range("N116").Value = -5601

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 750.5
$ws.Range("I3").Value = 663
$ws.Range("J3").Value = 1013
$ws.Range("K3").Value = 663
$ws.Range("L3").Value = 1013
$ws.Range("M3").Value = -549
$ws.Range("N3").Value = -1241

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2369.75
$ws.Range("I86").Value = 1833
$ws.Range("J86").Value = 3980
$ws.Range("K86").Value = 1833
$ws.Range("L86").Value = 3980
$ws.Range("M86").Value = -710
$ws.Range("N86").Value = -6226

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2369.75
$ws.Range("I89").Value = 1833
$ws.Range("J89").Value = 3980
$ws.Range("K89").Value = 9165
$ws.Range("L89").Value = 19900
$ws.Range("M89").Value = -3549
$ws.Range("N89").Value = -31132

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5429.1665
$ws.Range("I31").Value = 2805.5
$ws.Range("J31").Value = 6741
$ws.Range("K31").Value = 2805.5
$ws.Range("L31").Value = 6741
$ws.Range("M31").Value = -2510.5
$ws.Range("N31").Value = -7331

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5429.1665
$ws.Range("I34").Value = 2805.5
$ws.Range("J34").Value = 6741
$ws.Range("K34").Value = 2805.5
$ws.Range("L34").Value = 6741
$ws.Range("M34").Value = -2603.5
$ws.Range("N34").Value = -7145

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 10702.5
$ws.Range("I99").Value = 7810
$ws.Range("J99").Value = 11666.667
$ws.Range("K99").Value = 7810
$ws.Range("L99").Value = 11666.667
$ws.Range("M99").Value = -6312
$ws.Range("N99").Value = -14662.667

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 10702.5
$ws.Range("I126").Value = 7810
$ws.Range("J126").Value = 11666.667
$ws.Range("K126").Value = 23430
$ws.Range("L126").Value = 35000.001
$ws.Range("M126").Value = -20960
$ws.Range("N126").Value = -39940.001

# CRP row 139
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 48485
$ws.Range("J139").Value = 48485
$ws.Range("L139").Value = 48485
$ws.Range("N139").Value = -58765

# CUL row 12
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 47.782608
$ws.Range("J12").Value = 28.3125
$ws.Range("L12").Value = 84.9375
$ws.Range("N12").Value = -430.9375

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 952.4761999999999
$ws.Range("I113").Value = 736.9474
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 2210.8422
$ws.Range("L113").Value = 9000
$ws.Range("M113").Value = -40.84220000000005
$ws.Range("N113").Value = -13340

# GSM row 12
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").Value = ""

# GSM row 14
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 15005589
$ws.Range("I14").Value = 16875038
$ws.Range("J14").Value = 50000
$ws.Range("K14").Value = 16875038
$ws.Range("L14").Value = 50000
$ws.Range("M14").Value = -16874870
$ws.Range("N14").Value = -50336

# LTW row 20
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 2500
$ws.Range("J20").Value = 2500
$ws.Range("L20").Value = 2500
$ws.Range("N20").Value = -2952

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2752.6843
$ws.Range("J22").Value = 2577
$ws.Range("L22").Value = 2577
$ws.Range("N22").Value = -3167

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 2752.6843
$ws.Range("J27").Value = 2577
$ws.Range("L27").Value = 2577
$ws.Range("N27").Value = -2791

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6197.724
$ws.Range("I40").Value = 5718.9165
$ws.Range("J40").Value = 8496
$ws.Range("K40").Value = 5718.9165
$ws.Range("L40").Value = 8496
$ws.Range("M40").Value = -5582.9165
$ws.Range("N40").Value = -8768

# WVR row 23
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 57673.832
$ws.Range("I23").Value = 15336.667
$ws.Range("K23").Value = 15336.667
$ws.Range("M23").Value = -15107.667

# WVR row 47
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 35933
$ws.Range("J47").Value = 35933
$ws.Range("L47").Value = 35933
$ws.Range("N47").Value = -37077

# WVR row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 16622.223
$ws.Range("J54").Value = 16622.223
$ws.Range("L54").Value = 16622.223
$ws.Range("N54").Value = -17662.223

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 720514.2
$ws.Range("I126").Value = 4146.3335
$ws.Range("K126").Value = 12439.0005
$ws.Range("M126").Value = -9969.000499999998

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 55958.74
$ws.Range("J135").Value = 55958.74
$ws.Range("L135").Value = 55958.74
$ws.Range("N135").Value = -66098.73999999999
